$wb = $excel.ActiveWorkbook

# --- DataTable_Character: rename "armor" column to "accessori" and populate
#     helmet / accessori / weapon drop-item columns with real indices ------
$wsChar = $wb.Worksheets.Item("DataTable_Character")

# Header K1: "armor" -> "accessori"
$wsChar.Range("K1").Value = "accessori"

# J (helmet), K (accessori), L (weapon) drop item indices per row
$wsChar.Range("J2").Value = 14
$wsChar.Range("K2").Value = 13
$wsChar.Range("L2").Value = 15

$wsChar.Range("J3").Value = 25
$wsChar.Range("K3").Value = 24
$wsChar.Range("L3").Value = 26

$wsChar.Range("J4").Value = 19
$wsChar.Range("K4").Value = 18
$wsChar.Range("L4").Value = 20

$wsChar.Range("J5").Value = 3
$wsChar.Range("K5").Value = 21
$wsChar.Range("L5").Value = 4

$wsChar.Range("J6").Value = 11
$wsChar.Range("K6").Value = 10
$wsChar.Range("L6").Value = 12

$wsChar.Range("J7").Value = 16
$wsChar.Range("L7").Value = 17

$wsChar.Range("J8").Value = 8
$wsChar.Range("K8").Value = 7
$wsChar.Range("L8").Value = 9

$wsChar.Range("J9").Value = 5
$wsChar.Range("L9").Value = 6

$wsChar.Range("J10").Value = 22
$wsChar.Range("K10").Value = 21
$wsChar.Range("L10").Value = 23

# Bonus row: level 1 -> 10
$wsChar.Range("Q11").Value = 10

# --- sheet view / selection bookkeeping -----------------------------------
# DataTable_Item loses the active tab, selection moves to H19
$wsItem = $wb.Worksheets.Item("DataTable_Item")
$wsItem.Activate()
$wsItem.Range("H19").Select()

# DataTable_Character becomes the active tab, selection moves to J9
$wsChar.Activate()
$wsChar.Range("J9").Select()
